# Added duplicate insertion prevention:
# append new contact rows (each with its own "Name" mailto hyperlink) and
# populate a helper "some" index column used to flag/prevent duplicates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")

function Add-Contact($row, $name, $phone, $email) {
    $nameCell = $ws.Cells.Item($row, 1)
    $nameCell.Value = $name
    $nameCell.Hyperlinks.Add($nameCell, "mailto:" + $name)
    $nameCell.Style = "Hyperlink"

    $ws.Cells.Item($row, 2).Value = $phone
    $ws.Cells.Item($row, 3).Value = $email
}

# --- Append new contact rows (Name/Phone/Email) ------------------------
Add-Contact 7 "ced19i001@iiitdm.ac.in" 2 876789878
Add-Contact 8 "ced19i034@iiitdm.ac.in" 4 2123232444

# --- New "some" index column (D) ---------------------------------------
$ws.Cells.Item(1, 4).Value = "some"
$dValues = @(1, 2, 3, 4, 5, 5, 7)
for ($i = 0; $i -lt $dValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $dValues[$i]
}

# --- Last new contact row ------------------------------------------------
Add-Contact 9 "ced19i098@iiitdm.ac.in" 4 2123232444
$ws.Cells.Item(9, 4).Value = 7

# --- Restore the view's active selection -------------------------------
$ws.Range("E13").Select()
